$d = $word.ActiveDocument

$replacements = @(
    @("350÷3=", "956÷2="),
    @("298÷8=", "659÷6="),
    @("544÷8=", "758÷5="),
    @("579÷9=", "555÷3="),
    @("998÷6=", "783÷4="),
    @("486÷5=", "860÷5="),
    @("725÷3=", "108÷4="),
    @("488÷5=", "350÷4="),
    @("551÷7=", "724÷3="),
    @("581÷4=", "749÷6="),
    @("866÷2=", "657÷8="),
    @("831÷2=", "669÷9="),
    @("656÷6=", "327÷5="),
    @("810÷6=", "969÷7="),
    @("177÷9=", "437÷8="),
    @("387÷9=", "146÷4="),
    @("680÷4=", "183÷9="),
    @("399÷5=", "219÷8="),
    @("427÷4=", "222÷9="),
    @("211÷5=", "170÷8="),
    @("971÷8=", "343÷2="),
    @("794÷9=", "429÷9="),
    @("859÷7=", "236÷7="),
    @("195÷6=", "842÷6="),
    @("391÷4=", "132÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
